$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Requisitos:" label row (mirrors the A12/A17 section-header rows)
$ws.Range("A22").Value = "Requisitos:"

# New requirement text, duplicated into B23/C23 (col C carries the red-highlight style)
$text = "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)`n"
$ws.Range("B23").Value = $text
$ws.Range("C23").Value = $text

# Pick up the existing column B/C cell formatting (styles 2/3) from row 21
# instead of letting the new cells mint fresh style entries.
$ws.Range("B21:C21").Copy()
$ws.Range("B23:C23").PasteSpecial(-4122)

# Match the target row height for the new wrapped-text row.
$ws.Rows.Item(23).RowHeight = 30
